# --------------------------------------------------------------------------
# Fix PO analysis workbook:
#   1. Rename "Requested quantity" header to "Weekly_PO_Qty" on "Weekly Quantity"
#   2. Rename "Requested quantity" header to "Monthly_PO_Qty" on "Monthly Trend"
#   3. Add a new "PO Forecast" worksheet (ds / PO_Forecast / yhat_lower / yhat_upper)
# --------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Step 1: Update "Weekly Quantity" header text ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- Step 2: Update "Monthly Trend" header text ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Step 3: Create the new "PO Forecast" worksheet after the last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$newSheet.Name = "PO Forecast"

# Header row values
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Match the bold/centered/bordered header style used on the other sheets
$wsWeekly.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$newSheet.Cells.Item(2,1).Value = 45270.99999999999
$newSheet.Cells.Item(2,2).Value = 72
$newSheet.Cells.Item(2,3).Value = 9.083675389752933
$newSheet.Cells.Item(2,4).Value = 136.7819448314093
$newSheet.Cells.Item(3,1).Value = 45298.99999999999
$newSheet.Cells.Item(3,2).Value = 71
$newSheet.Cells.Item(3,3).Value = 4.396111160019692
$newSheet.Cells.Item(3,4).Value = 136.1324209959311
$newSheet.Cells.Item(4,1).Value = 45305.99999999999
$newSheet.Cells.Item(4,2).Value = 70
$newSheet.Cells.Item(4,3).Value = 3.31504815478847
$newSheet.Cells.Item(4,4).Value = 139.1745716882328
$newSheet.Cells.Item(5,1).Value = 45312.99999999999
$newSheet.Cells.Item(5,2).Value = 70
$newSheet.Cells.Item(5,3).Value = 3.694177918289841
$newSheet.Cells.Item(5,4).Value = 133.366690142545
$newSheet.Cells.Item(6,1).Value = 45319.99999999999
$newSheet.Cells.Item(6,2).Value = 70
$newSheet.Cells.Item(6,3).Value = 4.796762316345282
$newSheet.Cells.Item(6,4).Value = 130.0863788273732
$newSheet.Cells.Item(7,1).Value = 45333.99999999999
$newSheet.Cells.Item(7,2).Value = 69
$newSheet.Cells.Item(7,3).Value = 0.8115464121903048
$newSheet.Cells.Item(7,4).Value = 129.9361616740409
$newSheet.Cells.Item(8,1).Value = 45347.99999999999
$newSheet.Cells.Item(8,2).Value = 69
$newSheet.Cells.Item(8,3).Value = 9.463643291314131
$newSheet.Cells.Item(8,4).Value = 135.9199989862379
$newSheet.Cells.Item(9,1).Value = 45354.99999999999
$newSheet.Cells.Item(9,2).Value = 68
$newSheet.Cells.Item(9,3).Value = 2.690499092725108
$newSheet.Cells.Item(9,4).Value = 133.0486433295623
$newSheet.Cells.Item(10,1).Value = 45361.99999999999
$newSheet.Cells.Item(10,2).Value = 68
$newSheet.Cells.Item(10,3).Value = 3.377166599412183
$newSheet.Cells.Item(10,4).Value = 129.2633425177958
$newSheet.Cells.Item(11,1).Value = 45375.99999999999
$newSheet.Cells.Item(11,2).Value = 67
$newSheet.Cells.Item(11,3).Value = 0.6102292393870432
$newSheet.Cells.Item(11,4).Value = 133.42149215387
$newSheet.Cells.Item(12,1).Value = 45382.99999999999
$newSheet.Cells.Item(12,2).Value = 67
$newSheet.Cells.Item(12,3).Value = 7.262867209376824
$newSheet.Cells.Item(12,4).Value = 127.2474478363953
$newSheet.Cells.Item(13,1).Value = 45389.99999999999
$newSheet.Cells.Item(13,2).Value = 67
$newSheet.Cells.Item(13,3).Value = 6.355796383683987
$newSheet.Cells.Item(13,4).Value = 135.8774785036894
$newSheet.Cells.Item(14,1).Value = 45417.99999999999
$newSheet.Cells.Item(14,2).Value = 66
$newSheet.Cells.Item(14,3).Value = 4.017219189966046
$newSheet.Cells.Item(14,4).Value = 132.8108115133362
$newSheet.Cells.Item(15,1).Value = 45424.99999999999
$newSheet.Cells.Item(15,2).Value = 65
$newSheet.Cells.Item(15,3).Value = 5.496157418110355
$newSheet.Cells.Item(15,4).Value = 132.9246316347675
$newSheet.Cells.Item(16,1).Value = 45431.99999999999
$newSheet.Cells.Item(16,2).Value = 65
$newSheet.Cells.Item(16,3).Value = -4.042781374883705
$newSheet.Cells.Item(16,4).Value = 126.7482962598011
$newSheet.Cells.Item(17,1).Value = 45438.99999999999
$newSheet.Cells.Item(17,2).Value = 65
$newSheet.Cells.Item(17,3).Value = 0.1328274703126939
$newSheet.Cells.Item(17,4).Value = 127.2114694543218
$newSheet.Cells.Item(18,1).Value = 45459.99999999999
$newSheet.Cells.Item(18,2).Value = 64
$newSheet.Cells.Item(18,3).Value = 2.083291553160853
$newSheet.Cells.Item(18,4).Value = 127.6200654744113
$newSheet.Cells.Item(19,1).Value = 45466.99999999999
$newSheet.Cells.Item(19,2).Value = 64
$newSheet.Cells.Item(19,3).Value = 3.596178125891023
$newSheet.Cells.Item(19,4).Value = 128.6824309442698
$newSheet.Cells.Item(20,1).Value = 45473.99999999999
$newSheet.Cells.Item(20,2).Value = 63
$newSheet.Cells.Item(20,3).Value = 0.5183367310416106
$newSheet.Cells.Item(20,4).Value = 126.1715777968234
$newSheet.Cells.Item(21,1).Value = 45487.99999999999
$newSheet.Cells.Item(21,2).Value = 63
$newSheet.Cells.Item(21,3).Value = -6.121762628806001
$newSheet.Cells.Item(21,4).Value = 129.122501132057
$newSheet.Cells.Item(22,1).Value = 45508.99999999999
$newSheet.Cells.Item(22,2).Value = 62
$newSheet.Cells.Item(22,3).Value = 2.004122518131218
$newSheet.Cells.Item(22,4).Value = 122.8657344168548
$newSheet.Cells.Item(23,1).Value = 45522.99999999999
$newSheet.Cells.Item(23,2).Value = 62
$newSheet.Cells.Item(23,3).Value = 0.3245645936231329
$newSheet.Cells.Item(23,4).Value = 123.1856300554967
$newSheet.Cells.Item(24,1).Value = 45529.99999999999
$newSheet.Cells.Item(24,2).Value = 61
$newSheet.Cells.Item(24,3).Value = -4.527194677722289
$newSheet.Cells.Item(24,4).Value = 130.7903231470366
$newSheet.Cells.Item(25,1).Value = 45536.99999999999
$newSheet.Cells.Item(25,2).Value = 61
$newSheet.Cells.Item(25,3).Value = -4.261826676725969
$newSheet.Cells.Item(25,4).Value = 116.4243023631703
$newSheet.Cells.Item(26,1).Value = 45543.99999999999
$newSheet.Cells.Item(26,2).Value = 61
$newSheet.Cells.Item(26,3).Value = -3.424143204139174
$newSheet.Cells.Item(26,4).Value = 122.0071479139333
$newSheet.Cells.Item(27,1).Value = 45550.99999999999
$newSheet.Cells.Item(27,2).Value = 60
$newSheet.Cells.Item(27,3).Value = -4.666587861441401
$newSheet.Cells.Item(27,4).Value = 123.2042321955096
$newSheet.Cells.Item(28,1).Value = 45557.99999999999
$newSheet.Cells.Item(28,2).Value = 60
$newSheet.Cells.Item(28,3).Value = -1.116646256360271
$newSheet.Cells.Item(28,4).Value = 124.1016915502816
$newSheet.Cells.Item(29,1).Value = 45571.99999999999
$newSheet.Cells.Item(29,2).Value = 60
$newSheet.Cells.Item(29,3).Value = -2.573059494893128
$newSheet.Cells.Item(29,4).Value = 122.215142255188
$newSheet.Cells.Item(30,1).Value = 45578.99999999999
$newSheet.Cells.Item(30,2).Value = 59
$newSheet.Cells.Item(30,3).Value = -3.998190352721078
$newSheet.Cells.Item(30,4).Value = 126.6729138048099
$newSheet.Cells.Item(31,1).Value = 45592.99999999999
$newSheet.Cells.Item(31,2).Value = 59
$newSheet.Cells.Item(31,3).Value = -5.620764583155692
$newSheet.Cells.Item(31,4).Value = 121.2722966148737
$newSheet.Cells.Item(32,1).Value = 45599.99999999999
$newSheet.Cells.Item(32,2).Value = 58
$newSheet.Cells.Item(32,3).Value = -7.988981466917333
$newSheet.Cells.Item(32,4).Value = 122.4343683238102
$newSheet.Cells.Item(33,1).Value = 45606.99999999999
$newSheet.Cells.Item(33,2).Value = 58
$newSheet.Cells.Item(33,3).Value = -9.264757561239932
$newSheet.Cells.Item(33,4).Value = 125.7078373940284
$newSheet.Cells.Item(34,1).Value = 45613.99999999999
$newSheet.Cells.Item(34,2).Value = 58
$newSheet.Cells.Item(34,3).Value = -2.971953756440724
$newSheet.Cells.Item(34,4).Value = 120.8874252159201
$newSheet.Cells.Item(35,1).Value = 45620.99999999999
$newSheet.Cells.Item(35,2).Value = 58
$newSheet.Cells.Item(35,3).Value = -3.218637580814448
$newSheet.Cells.Item(35,4).Value = 119.1618115917357
$newSheet.Cells.Item(36,1).Value = 45627.99999999999
$newSheet.Cells.Item(36,2).Value = 57
$newSheet.Cells.Item(36,3).Value = -8.791146593422026
$newSheet.Cells.Item(36,4).Value = 118.2581069895916
$newSheet.Cells.Item(37,1).Value = 45634.99999999999
$newSheet.Cells.Item(37,2).Value = 57
$newSheet.Cells.Item(37,3).Value = -3.787724267059895
$newSheet.Cells.Item(37,4).Value = 119.3642933056399
$newSheet.Cells.Item(38,1).Value = 45641.99999999999
$newSheet.Cells.Item(38,2).Value = 57
$newSheet.Cells.Item(38,3).Value = -3.292267590225886
$newSheet.Cells.Item(38,4).Value = 119.6316606910875
$newSheet.Cells.Item(39,1).Value = 45648.99999999999
$newSheet.Cells.Item(39,2).Value = 56
$newSheet.Cells.Item(39,3).Value = -11.37089370295639
$newSheet.Cells.Item(39,4).Value = 121.7546522501072
$newSheet.Cells.Item(40,1).Value = 45655.99999999999
$newSheet.Cells.Item(40,2).Value = 56
$newSheet.Cells.Item(40,3).Value = -9.247280591488266
$newSheet.Cells.Item(40,4).Value = 120.9912476430019

# Match the date-time number format used for the "ds" column on the other sheets
$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A40").PasteSpecial(-4122)

$excel.CutCopyMode = 0

Write-Output "PO Forecast sheet created with $($newSheet.UsedRange.Rows.Count) rows"
